# "Fix Resumo de publicados"
# Updates status/edit-log cells on several existing rows and appends six
# new credit records (rows 221-226) to the "base_filtrada" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column map (header row):
# A Situação | B Origem de Recursos | C Órgão (UO) | D Nº do Processo
# E Tipo de Crédito | F Fonte de Recursos | G Grupo de Despesas | H Valor
# I Objetivo | J Observação | K Data de Recebimento | L Data de Publicação
# M Nº do decreto | N Contabilizar no Limite? | O Cadastrado Por
# P Última Edição | Q Nº ATA | R Opnião SOP

function Set-DateCell($row, $col, $serial) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $serial
    $c.NumberFormat = "YYYY-MM-DD"
}

# ---------------------------------------------------------------------
# Updates to existing rows
# ---------------------------------------------------------------------

# Row 177
$ws.Cells.Item(177, 2).Value = "Sem Cobertura - Atendido pela 1º NR"
$ws.Cells.Item(177, 16).Value = "Fillipe Fernando - 28/05/2025 18:07:50"

# Row 186
$ws.Cells.Item(186, 1).Value = "Análise - SEFAZ"
$ws.Cells.Item(186, 16).Value = "Romero Cavalcanti - 29/05/2025 14:04:46"

# Row 194
$ws.Cells.Item(194, 2).Value = "Sem Cobertura - Atendido pela 1º NR"
$ws.Cells.Item(194, 16).Value = "Fillipe Fernando - 28/05/2025 18:09:34"

# Row 212
$ws.Cells.Item(212, 1).Value = "BLOCO 434050 - SOP - Superintendente"
$ws.Cells.Item(212, 16).Value = "Isadora Costa - 29/05/2025 17:51:26"

# Row 214
$ws.Cells.Item(214, 1).Value = "Publicado"
$ws.Cells.Item(214, 2).Value = "Sem Cobertura - Atendido pela 1º NR"
Set-DateCell 214 12 45806
$ws.Cells.Item(214, 13).Value = 102469
$ws.Cells.Item(214, 16).Value = "Felliphy Queiroz - 29/05/2025 12:01:01"

# Row 216
$ws.Cells.Item(216, 2).Value = "Redução/Anulação"
$ws.Cells.Item(216, 16).Value = "Isadora Costa - 29/05/2025 17:35:34"

# Row 217
$ws.Cells.Item(217, 1).Value = "BLOCO 434066 - SEPLAG - Demais Orgãos"
$ws.Cells.Item(217, 16).Value = "Isadora Costa - 29/05/2025 18:08:07"

# Row 218
$ws.Cells.Item(218, 1).Value = "BLOCO 434066 - SEPLAG - Demais Orgãos"
$ws.Cells.Item(218, 16).Value = "Isadora Costa - 29/05/2025 17:18:33"

# ---------------------------------------------------------------------
# New rows 221-226
# ---------------------------------------------------------------------

$newRows = @(
    @{ Row=221; A="Análise - SOP"; B="Redução/Anulação"; C="PMAL";
       D="E:01206.0000032071/2025"; E="Suplementar"; F=500; G=3; H=7838000;
       I="Complemento orçamentário para cumprir com as demandas inerentes a atividade Policial.";
       J=""; K=45805; N="SIM";
       O="Romero Cavalcanti - 28/05/2025 17:45:22"; P="" },

    @{ Row=222; A="BLOCO 434078 - SEFAZ - Despachos e Decretos"; B="Sem Cobertura - Atendido pela 1º NR"; C="SECOM";
       D="E:01700.0000003990/2025"; E="Suplementar"; F=500; G=3; H=8155000;
       I="Atender a execução das atividades de publicidade planejadas para o primeiro bimestre de 2025.";
       J="ata 26"; K=45805; N="SIM";
       O="Fillipe Fernando - 28/05/2025 18:03:33"; P="Fillipe Fernando - 28/05/2025 20:48:39" },

    @{ Row=223; A="Análise - CPOF"; B="Sem Cobertura"; C="SECOM";
       D="E:02200.0000000671/2025"; E="Suplementar"; F=500; G=3; H=12700000;
       I="Atender a ações previstas para os meses de junho e julho do 2025.";
       J=""; K=45806; N="SIM";
       O="Fillipe Fernando - 29/05/2025 13:12:08"; P="Romero Cavalcanti - 29/05/2025 17:56:26" },

    @{ Row=224; A="Análise - SOP"; B="Sem Cobertura"; C="SETRAND";
       D="E:35032.0000001392/2025"; E="Suplementar"; F=754; G=4; H=8313500;
       I="DUPLICAÇÃO, RESTAURAÇÃO COM MELHORIAS DA RODOVIA AL 101 NORTE, TRECHO: MACEIÓ/BARRA DE SANTO ANTÔNIO";
       J=""; K=45807; N="SIM";
       O="Fillipe Fernando - 30/05/2025 16:00:12"; P="" },

    @{ Row=225; A="Análise - SOP"; B="Redução/Anulação"; C="FUNTURIS";
       D="E:29032.0000000514/2025"; E="Suplementar"; F=759; G=3; H=400000;
       I="Ampliação e consolidação do destino Alagoas nos mercados Nacional e Internacional do Fundo do Turismo – FUNTURIS da Secretaria de Estado do Turismo - SETUR.";
       J=""; K=45807; N="SIM";
       O="Fillipe Fernando - 30/05/2025 16:03:58"; P="" },

    @{ Row=226; A="Análise - SOP"; B="Sem Cobertura"; C="SETUR";
       D="E:29032.0000000509/2025"; E="Suplementar"; F=500; G=3; H=4813155.82;
       I="Atender a diversos contratos de manutenção e eventos.";
       J=""; K=45807; N="SIM";
       O="Fillipe Fernando - 30/05/2025 16:07:28"; P="" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    Set-DateCell $row 11 $r.K
    $ws.Cells.Item($row, 12).Value = ""
    $ws.Cells.Item($row, 13).Value = ""
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = ""
    $ws.Cells.Item($row, 18).Value = ""
}
